$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/10/2025  Through  2/16/2025"

# --- Crime-complaints grid updates ---
# Row 14
$ws.Range("D14").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = "'***.*"
$ws.Range("C15").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("G14").Value = 3
# Row 15
$ws.Range("D15").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -62.5
$ws.Range("N15").Value = -62.5
# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 15
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 42
$ws.Range("K16").Value = -38.095238095238
$ws.Range("L16").Value = -16.129032258064
$ws.Range("M16").Value = -39.53488372093
$ws.Range("N16").Value = -91.216216216216
# Row 17
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 43
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = 10.25641025641
$ws.Range("I17").Value = 77
$ws.Range("J17").Value = 73
$ws.Range("K17").Value = 5.479452054794
$ws.Range("L17").Value = -14.444444444444
$ws.Range("M17").Value = 28.333333333333
$ws.Range("N17").Value = -38.4
# Row 18
$ws.Range("C18").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 6
$ws.Range("H18").Value = -33.333333333333
$ws.Range("J18").Value = 15
$ws.Range("K18").Value = -33.333333333333
$ws.Range("L18").Value = -61.538461538461
$ws.Range("M18").Value = -72.972972972973
$ws.Range("N18").Value = -88.372093023255
# Row 19
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -37.5
$ws.Range("F19").Value = 14
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -56.25
$ws.Range("I19").Value = 26
$ws.Range("J19").Value = 54
$ws.Range("K19").Value = -51.851851851851
$ws.Range("L19").Value = -48
$ws.Range("M19").Value = -43.478260869565
$ws.Range("N19").Value = -69.767441860465
# Row 20
$ws.Range("C20").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -20
$ws.Range("J20").Value = 24
$ws.Range("K20").Value = -41.666666666666
$ws.Range("L20").Value = -6.666666666666
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -84.782608695652
# Row 21
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -22.727272727272
$ws.Range("F21").Value = 87
$ws.Range("G21").Value = 126
$ws.Range("H21").Value = -30.952380952381
$ws.Range("I21").Value = 156
$ws.Range("J21").Value = 215
$ws.Range("K21").Value = -27.441860465116
$ws.Range("L21").Value = -29.729729729729
$ws.Range("M21").Value = -24.271844660194
$ws.Range("N21").Value = -77.74607703281
# Row 22
$ws.Range("C22").Value = 2
$ws.Range("D16").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("C15").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -77.777777777777
$ws.Range("I22").Value = 5
$ws.Range("K22").Value = -68.75
$ws.Range("L22").Value = -44.444444444444
$ws.Range("M22").Value = 66.666666666666
# Row 23
$ws.Range("C23").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 31
$ws.Range("H23").Value = -58.064516129032
$ws.Range("J23").Value = 48
$ws.Range("K23").Value = -35.416666666666
$ws.Range("L23").Value = -38
$ws.Range("M23").Value = 19.230769230769
# Row 24
$ws.Range("C24").Value = 15
$ws.Range("E24").Value = -48.275862068965
$ws.Range("F24").Value = 73
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = -32.407407407407
$ws.Range("I24").Value = 111
$ws.Range("J24").Value = 162
$ws.Range("K24").Value = -31.481481481481
$ws.Range("L24").Value = -30.625
$ws.Range("M24").Value = -6.72268907563
# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -72.727272727272
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = -59.375
$ws.Range("I25").Value = 22
$ws.Range("J25").Value = 48
$ws.Range("K25").Value = -54.166666666666
$ws.Range("L25").Value = -47.619047619047
# Row 26
$ws.Range("C26").Value = 22
$ws.Range("E26").Value = 29.411764705882
$ws.Range("F26").Value = 67
$ws.Range("G26").Value = 63
$ws.Range("H26").Value = 6.349206349206
$ws.Range("I26").Value = 114
$ws.Range("J26").Value = 99
$ws.Range("K26").Value = 15.151515151515
$ws.Range("L26").Value = -5.785123966942
$ws.Range("M26").Value = -26.451612903225
# Row 27
$ws.Range("D27").Value = 1
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = -40
$ws.Range("L27").Value = -66.666666666666
# Row 28
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 1
$ws.Range("D16").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = 200
$ws.Range("E16").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 8
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 15
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = 87.5
$ws.Range("L28").Value = 15.384615384615
# Row 29
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 3
$ws.Range("K29").Value = -62.5
$ws.Range("L29").Value = -50
$ws.Range("M29").Value = -62.5
$ws.Range("N29").Value = -88.461538461538
# Row 30
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 3
$ws.Range("K30").Value = -57.142857142857
$ws.Range("L30").Value = -50
$ws.Range("N30").Value = -88.461538461538
